$d = $word.ActiveDocument

# Find the paragraph that contains the sentence to be removed
# ("und dann starben sie ") and delete the whole paragraph,
# including its trailing paragraph mark, so the preceding
# paragraph (ending in a single space) becomes directly followed
# by the page-break paragraph again.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*und dann starben sie*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
